$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D that are purely numeric-looking text (e.g. "63.499.60")
# are prefixed with a literal leading apostrophe so Excel keeps them as text
# (matching the source data's "thousands-dot" price formatting) instead of
# auto-converting to a number and silently dropping significant trailing zeros.

$ws.Range('D2').Value = "'" + '63.499.60'
$ws.Range('E2').Value = '  +0.99%  '
$ws.Range('D3').Value = "'" + '3.099.42'
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = "'" + '582.91'
$ws.Range('E5').Value = '  -0.12%  '
$ws.Range('D6').Value = "'" + '144.78'
$ws.Range('E6').Value = '  -0.19%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').Value = "'" + '3.091.85'
$ws.Range('E8').Value = '  -0.47%  '
$ws.Range('D9').Value = "'" + '0.527'
$ws.Range('E9').Value = '  -0.41%  '
$ws.Range('E10').Value = '  +6.61%  '
$ws.Range('D11').Value = "'" + '5.60'
$ws.Range('E11').Value = '  -2.81%  '
$ws.Range('E12').Value = '  -2.56%  '
$ws.Range('D13').Value = "'" + '0.0000245'
$ws.Range('E13').Value = '  -0.91%  '
$ws.Range('D14').Value = "'" + '37.15'
$ws.Range('E14').Value = '  +4.33%  '
$ws.Range('E15').Value = '  -1.14%  '
$ws.Range('D16').Value = "'" + '3.612.42'
$ws.Range('E16').Value = '  -0.49%  '
$ws.Range('D17').Value = "'" + '63.356.91'
$ws.Range('E17').Value = '  +0.91%  '
$ws.Range('E18').Value = '  -1.00%  '
$ws.Range('D19').Value = "'" + '3.095.07'
$ws.Range('E19').Value = '  -0.45%  '
$ws.Range('D20').Value = "'" + '461.70'
$ws.Range('E20').Value = '  -1.25%  '
$ws.Range('D21').Value = "'" + '14.22'
$ws.Range('E21').Value = '  +1.03%  '
$ws.Range('D22').Value = "'" + '0.724'
$ws.Range('E22').Value = '  -0.61%  '
$ws.Range('D23').Value = "'" + '7.45'
$ws.Range('E23').Value = '  -1.26%  '
$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').Value = "'" + '81.31'
$ws.Range('E24').Value = '  -0.87%  '
$ws.Range('B25').Value = 'InternetComputer(DFINITY)'
$ws.Range('C25').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D25').Value = "'" + '12.92'
$ws.Range('E25').Value = '  -3.01%  '
$ws.Range('E26').Value = '  -1.60%  '
$ws.Range('D28').Value = "'" + '9.00'
$ws.Range('E28').Value = '  +8.97%  '
$ws.Range('E29').Value = '  +0.06%  '
$ws.Range('E30').Value = '  -0.46%  '
$ws.Range('D31').Value = "'" + '2.19'
$ws.Range('E31').Value = '  -1.90%  '
$ws.Range('D32').Value = "'" + '6.83'
$ws.Range('E32').Value = '  +0.24%  '
$ws.Range('D33').Value = "'" + '0.110'
$ws.Range('E33').Value = '  -0.96%  '
$ws.Range('D34').Value = "'" + '26.62'
$ws.Range('E34').Value = '  -1.40%  '
$ws.Range('D35').Value = '0.0₃0851'
$ws.Range('E35').Value = '  -2.17%  '
$ws.Range('E36').Value = '  +3.68%  '
$ws.Range('B37').Value = 'Mantle'
$ws.Range('C37').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D37').Value = "'" + '1.03'
$ws.Range('E37').Value = '  -1.14%  '
$ws.Range('B38').Value = 'Stacks'
$ws.Range('C38').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D38').Value = "'" + '2.30'
$ws.Range('E38').Value = '  -3.20%  '
$ws.Range('E39').Value = '  -1.02%  '
$ws.Range('D40').Value = "'" + '50.24'
$ws.Range('E40').Value = '  -1.37%  '
$ws.Range('D41').Value = "'" + '434.71'
$ws.Range('E41').Value = '  +0.10%  '
$ws.Range('D42').Value = "'" + '8.72'
$ws.Range('E42').Value = '  -0.26%  '
$ws.Range('E43').Value = '  -0.44%  '
$ws.Range('D44').Value = "'" + '2.876.30'
$ws.Range('E44').Value = '  -1.88%  '
$ws.Range('E45').Value = '  -3.26%  '
$ws.Range('E46').Value = '  -2.77%  '
$ws.Range('D47').Value = "'" + '35.75'
$ws.Range('E47').Value = '  +0.11%  '
$ws.Range('E48').Value = '  +0.02%  '
$ws.Range('D49').Value = "'" + '123.25'
$ws.Range('E49').Value = '  -0.30%  '
$ws.Range('E50').Value = '  -1.38%  '
$ws.Range('D51').Value = "'" + '24.08'
$ws.Range('E51').Value = '  -2.18%  '
